# daily auto push: 2025-10-02 07:26 UTC
# Append the new daily data row (row 51) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 51

# Columns A (date) and B (weekday) hold text that looks like a date, not an
# actual date value. Writing the literal string straight into .Value would
# make Excel auto-convert "2025/10/02" into a date serial and stamp the cell
# with a new number-format style. Instead we build the text via a formula
# (which always yields a string result) and then paste-special just the
# value back over itself; this keeps the cell's type as text/string without
# leaving any new style behind.
$ws.Range("A" + $row).Formula = "=""2025/10/02"""
$ws.Range("A" + $row).Copy() | Out-Null
$ws.Range("A" + $row).PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B" + $row).Value = "木"
$ws.Range("C" + $row).Value = 16
$ws.Range("D" + $row).Value = 26
